$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 4 header columns: Physical Hosts | Hostname |
# VirtualizationNode | CappedCPU. We are inserting a new
# "PhysicalServerModelName" header as the new column D, and pushing the
# existing "CappedCPU" header out to the new column E.

# 1. Clone column D (with its header formatting: bold, centered "Calibri"
#    style) into the new column E, so the moved header keeps the exact
#    same look as the other headers.
$ws.Range("D1").Copy($ws.Range("E1"))

# 2. Column E now holds the old "CappedCPU" header.
$ws.Range("E1").Value = "CappedCPU"

# 3. Column D becomes the new "PhysicalServerModelName" header (same
#    style/format it already had).
$ws.Range("D1").Value = "PhysicalServerModelName"

# 4. Widen column D to fit the new, longer header text, and give column E
#    the width that column D used to have (the width that fit "CappedCPU").
$ws.Columns.Item(4).ColumnWidth = 22.3333333333333
$ws.Columns.Item(5).ColumnWidth = 11.1666666666667

# 5. Reset the sheet's active selection back to A1.
$ws.Range("A1").Select() | Out-Null
